# Assignment 1 Presentation.pptx edit
#   1) The cached "today" text inside the auto-updating date field
#      (type="datetimeFigureOut") on the Slide Master and on every one
#      of the 11 Custom Layouts flips from 2/17/2020 -> 2/18/2020.
#   2) Slide 1's title textbox "TNC MCU Assignment 1" has its first two
#      words swapped -> "MCU TNC Assignment 1" (split across two runs).

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "2/17/2020") {
                $tr.Text = "2/18/2020"
            }
        }
    }
}

# Slide Master's own Date Placeholder shape.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every Custom Layout hanging off the master also carries its own
# (independent) Date Placeholder shape/field.
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Slide 1 title: "TNC MCU Assignment 1" -> "MCU TNC Assignment 1"
$s1 = $p.Slides.Item(1)
$title = $s1.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange

# Trim the run down to just the trailing "1" (keeps its dirty="0" rPr),
# then type the replacement lead-in text in front of it as a new run.
$lead = $titleRange.Characters(1, 19)
$lead.Text = "MCU TNC Assignment "
